# Update symbol list values as published by GitHub Actions crypto scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $CellRef, $Text)
    $cell = $Worksheet.Range($CellRef)
    # Force the cell to be treated as text so numeric-looking strings (prices)
    # are not silently converted into number cells, then restore a plain
    # "Normal" style so no extra number-format styling is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# --- Simple price (column D) updates ---
Set-TextValue $ws "D2"  "245.58"
Set-TextValue $ws "D3"  "24.11"
Set-TextValue $ws "D4"  "5.273"
Set-TextValue $ws "D5"  "0.05786"
Set-TextValue $ws "D6"  "6.476"
Set-TextValue $ws "D7"  "3.125"
Set-TextValue $ws "D8"  "0.8163"
Set-TextValue $ws "D9"  "0.8502"
Set-TextValue $ws "D11" "0.06924"
Set-TextValue $ws "D12" "0.03130"
Set-TextValue $ws "D13" "0.02893"
Set-TextValue $ws "D14" "0.09383"
Set-TextValue $ws "D15" "3.734"
Set-TextValue $ws "D16" "0.001534"
Set-TextValue $ws "D18" "0.0005961"
Set-TextValue $ws "D19" "0.006224"
Set-TextValue $ws "D20" "0.001233"
Set-TextValue $ws "D21" "0.004621"
Set-TextValue $ws "D22" "0.00006893"
Set-TextValue $ws "D23" "3.500"
Set-TextValue $ws "D24" "2.142"
Set-TextValue $ws "D25" "0.3193"
Set-TextValue $ws "D26" "0.1318"
Set-TextValue $ws "D28" "0.0002328"
Set-TextValue $ws "D40" "0.03652"

# --- Rows 41-43: coin ranking reshuffled (Kick/BKEX/CEJI rotate down one rank) ---
Set-TextValue $ws "B41" "KickToken"
Set-TextValue $ws "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D41" "0.006224"
Set-TextValue $ws "E41" "40KickTokenKICK"

Set-TextValue $ws "B42" "BKEXToken"
Set-TextValue $ws "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D42" "0.1053"
Set-TextValue $ws "E42" "41BKEXTokenBKK"

Set-TextValue $ws "B43" "CEJI"
Set-TextValue $ws "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D43" "0.002717"
Set-TextValue $ws "E43" "42CEJICEJI"

# --- Remaining price / label updates ---
Set-TextValue $ws "D44" "0.008389"
Set-TextValue $ws "D45" "0.00005246"
Set-TextValue $ws "D46" "0.00000000749"
Set-TextValue $ws "D47" "0.3694"
Set-TextValue $ws "E47" "46CoinbaseStockTokenCOINWorstin24h"
Set-TextValue $ws "D48" "0.002270"
Set-TextValue $ws "D49" "0.00002097"
Set-TextValue $ws "D50" "0.0001997"
